# Estadisticos Segundo Parcial 23 Mayo
#
# - Fills in the real "Segundo Parcial" (2P) figures on the "Estadisticos 2P"
#   sheet (previously just placeholder/duplicate numbers).
# - Recomputes the affected "Estadisticos Final" figures that depend on 2P.
# - Replaces/extends the "Rescatables" (make-up exam) roster with the
#   updated list of students (22 rows instead of 14).

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws4 = $wb.Worksheets.Item("Rescatables")

# ---------------------------------------------------------------------------
# Estadisticos 2P  -  columns: Mat, Grupo, Totales, Blancos, Reprobados,
#                               Aprobados, Por_Apro, Promedio
# Only D:H change (C "Totales" already correct); H ("Promedio") is new.
# ---------------------------------------------------------------------------
$data2 = @(
    @(0, 8, 29, 78.38, 7.8),
    @(0, 0, 39, 100, 8.5),
    @(0, 5, 27, 84.38, 7.8),
    @(0, 4, 35, 89.73999999999999, 8.1),
    @(0, 6, 32, 84.20999999999999, 8.1),
    @(0, 3, 21, 87.5, 8.5),
    @(0, 0, 34, 100, 8.9),
    @(0, 0, 31, 100, 9.1),
    @(0, 1, 39, 97.5, 9.199999999999999),
    @(0, 0, 24, 100, 8.199999999999999),
    @(0, 3, 35, 92.11, 8.9),
    @(0, 1, 27, 96.43000000000001, 8.6),
    @(0, 6, 25, 80.65000000000001, 8)
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $row = $i + 2
    $vals = $data2[$i]
    $ws2.Cells.Item($row, 4).Value = $vals[0]   # D Blancos
    $ws2.Cells.Item($row, 5).Value = $vals[1]   # E Reprobados
    $ws2.Cells.Item($row, 6).Value = $vals[2]   # F Aprobados
    $ws2.Cells.Item($row, 7).Value = $vals[3]   # G Por_Apro
    $ws2.Cells.Item($row, 8).Value = $vals[4]   # H Promedio
}

# ---------------------------------------------------------------------------
# Estadisticos Final - recompute the cumulative figures that shifted because
# of the new 2P numbers above. Only the affected cells are touched.
# ---------------------------------------------------------------------------
$ws3.Cells.Item(2, 5).Value = 8
$ws3.Cells.Item(2, 6).Value = 29
$ws3.Cells.Item(2, 7).Value = 78.38

$ws3.Cells.Item(3, 8).Value = 8.699999999999999

$ws3.Cells.Item(4, 5).Value = 5
$ws3.Cells.Item(4, 6).Value = 27
$ws3.Cells.Item(4, 7).Value = 84.38
$ws3.Cells.Item(4, 8).Value = 7.9

$ws3.Cells.Item(5, 5).Value = 4
$ws3.Cells.Item(5, 6).Value = 35
$ws3.Cells.Item(5, 7).Value = 89.73999999999999
$ws3.Cells.Item(5, 8).Value = 7.8

$ws3.Cells.Item(6, 5).Value = 6
$ws3.Cells.Item(6, 6).Value = 32
$ws3.Cells.Item(6, 7).Value = 84.20999999999999
$ws3.Cells.Item(6, 8).Value = 8

$ws3.Cells.Item(7, 8).Value = 8.1

$ws3.Cells.Item(8, 8).Value = 8.6

$ws3.Cells.Item(9, 8).Value = 8.800000000000001

$ws3.Cells.Item(10, 8).Value = 9.300000000000001

$ws3.Cells.Item(12, 5).Value = 3
$ws3.Cells.Item(12, 6).Value = 35
$ws3.Cells.Item(12, 7).Value = 92.11
$ws3.Cells.Item(12, 8).Value = 9

$ws3.Cells.Item(14, 8).Value = 7.8

# ---------------------------------------------------------------------------
# Rescatables - updated roster (22 students needing make-up work, sorted by
# number of failed subjects, descending). Columns: NC, Paterno, Materno,
# Nombres, Nombre_Largo (Mat), Grupo, Reprobadas.
# ---------------------------------------------------------------------------
$data4 = @(
    @(24330051920233, "ARELLANO",    "PAZ",        "ADRIAN",          "Ciencias sociales II",  "2AEM",  4),
    @(23330051920225, "FLORES",      "VAZQUEZ",    "MARCO ANTONIO",   "Ciencias sociales II",  "2AEM",  4),
    @(24330051920353, "GARCIA",      "SANCHEZ",    "JOY JARA",        "Ciencias sociales II",  "2AEM",  4),
    @(24330051920259, "RUIZ",        "PELLICO",    "YOSHUA RAFAEL",   "Ciencias sociales II",  "2APM",  4),
    @(24330051920404, "PARADA",      "SANTOS",     "MARCO DIDIEL",    "Ciencias sociales II",  "2BEM",  4),
    @(24330051920373, "MAZABA",      "QUINTERO",   "MANUEL ABDUL",    "Ciencias sociales II",  "2BLCM", 4),
    @(23330051920001, "ACEVEDO",     "GARCIA",     "OSCAR ORLANDO",   "Ciencias sociales II",  "2AEM",  3),
    @(24330051920387, "NICANOR",     "MALDONADO",  "DENISSE ARELI",   "Ciencias sociales II",  "2APM",  3),
    @(24330051920084, "VELAZQUEZ",   "LOPEZ",      "ADBIEL",          "Ciencias sociales II",  "2BEM",  3),
    @(22330051920010, "CRUZ",        "COYOHUA",    "DIEGO",           "TEMAS DE FILOSOFÍA",    "6AEM",  3),
    @(22330051920031, "CASTILLO",    "GONZALEZ",   "RICARDO",         "TEMAS DE FILOSOFÍA",    "6BEM",  3),
    @(22330051920189, "JENKINS",     "GARCIA",     "ARTHUR RICHARD",  "TEMAS DE FILOSOFÍA",    "6BEM",  3),
    @(22330051920043, "PALOMINO",    "HERNANDEZ",  "AARON MIGUEL",    "TEMAS DE FILOSOFÍA",    "6BEM",  3),
    @(24330051920103, "ESCOBAR",     "JUAN",       "GIOVANNI ARIEL",  "Ciencias sociales II",  "2AEM",  2),
    @(24330051920145, "HUERTA",      "GONZALEZ",   "YERIEL",          "Ciencias sociales II",  "2AEM",  2),
    @(24330051920109, "PERALTA",     "TRONCO",     "JONATHAN ISRAEL", "Ciencias sociales II",  "2BEM",  2),
    @(23330051920071, "RODRIGUEZ",   "CASTILLO",   "IVONNE ERIMAR",   "Ciencias sociales III", "4ARHM", 2),
    @(22330061460232, "ALVAREZ",     "VOTE",       "CAMILO",          "Ciencias sociales III", "4BLCM", 2),
    @(22330051920424, "COLMENARES",  "MARTINEZ",   "JULIO EDUARDO",   "Ciencias sociales III", "4BLCM", 2),
    @(22330051920033, "CRESCENCIO",  "DIAZ",       "DIEGO ARMANDO",   "TEMAS DE FILOSOFÍA",    "6BEM",  2),
    @(22330051920045, "PEREZ",       "ROMERO",     "JULIAN DAVID",    "TEMAS DE FILOSOFÍA",    "6BEM",  2)
)

for ($i = 0; $i -lt $data4.Length; $i++) {
    $row = $i + 2
    $vals = $data4[$i]
    $ws4.Cells.Item($row, 1).Value = $vals[0]   # A NC
    $ws4.Cells.Item($row, 2).Value = $vals[1]   # B Paterno
    $ws4.Cells.Item($row, 3).Value = $vals[2]   # C Materno
    $ws4.Cells.Item($row, 4).Value = $vals[3]   # D Nombres
    $ws4.Cells.Item($row, 5).Value = $vals[4]   # E Nombre_Largo (Mat)
    $ws4.Cells.Item($row, 6).Value = $vals[5]   # F Grupo
    $ws4.Cells.Item($row, 7).Value = $vals[6]   # G Reprobadas
}
